function Set-Row($ws, $r, $a, $b, $c, $d, $e) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet LP1912: header metadata ---
$ws1.Range("A2").Value = "Última actualización: 09:26:30"
$ws1.Range("A3").Value = "Total filas: 120"

# --- Sheet LP1912: data rows that changed (re-sort + newly scraped rows) ---
Set-Row $ws1 95 "07:24:45" "09:22" "17_ROMERO" 118 "LP1912"
Set-Row $ws1 96 "08:00:32" "09:22" "16_SANTA ANA" 82 "LP1912"
Set-Row $ws1 107 "09:26:30" "09:58" "16_SANTA ANA" 32 "LP1912"
Set-Row $ws1 108 "08:31:53" "10:03" "11_ETCHEVERRY" 92 "LP1912"
Set-Row $ws1 109 "09:26:30" "10:04" "23_HERNANDEZ" 38 "LP1912"
Set-Row $ws1 110 "08:31:53" "10:10" "16_P MOR-SANTA ANA" 99 "LP1912"
Set-Row $ws1 111 "08:31:53" "10:12" "15_ABASTO" 101 "LP1912"
Set-Row $ws1 112 "09:26:30" "10:13" "10_OLMOS" 47 "LP1912"
Set-Row $ws1 113 "08:31:53" "10:20" "26_HERNANDEZ" 109 "LP1912"
Set-Row $ws1 114 "08:55:25" "10:21" "26_HERNANDEZ" 86 "LP1912"
Set-Row $ws1 115 "08:31:53" "10:22" "17_ROMERO" 111 "LP1912"
Set-Row $ws1 116 "09:26:30" "10:23" "11_ETCHEVERRY" 57 "LP1912"
Set-Row $ws1 117 "08:31:53" "10:26" "215A_EL PATO" 115 "LP1912"
Set-Row $ws1 118 "08:47:51" "10:41" "17_ROMERO" 114 "LP1912"
Set-Row $ws1 119 "08:55:25" "10:42" "17_ROMERO" 107 "LP1912"
Set-Row $ws1 120 "08:47:51" "10:43" "14_ABASTO" 116 "LP1912"
Set-Row $ws1 121 "09:26:30" "10:57" "27_EL RETIRO" 91 "LP1912"
Set-Row $ws1 122 "09:26:30" "11:02" "215C_EL PATO" 96 "LP1912"
Set-Row $ws1 123 "09:26:30" "11:06" "16_P MOR-167 Y 521" 100 "LP1912"
Set-Row $ws1 124 "09:26:30" "11:19" "86_EST CHICA-ESC AGRARIA" 113 "LP1912"
Set-Row $ws1 125 "09:26:30" "11:21" "26_HERNANDEZ" 115 "LP1912"

# --- Sheet LP1912-215: header metadata ---
$ws2.Range("A2").Value = "Última actualización: 09:26:30"
$ws2.Range("A3").Value = "Total filas: 18"

# --- Sheet LP1912-215: new row appended ---
Set-Row $ws2 23 "09:26:30" "11:02" "215C_EL PATO" 96 "LP1912"

# --- Sheet 6203-6173: header metadata ---
$ws3.Range("A2").Value = "Última actualización: 09:26:30"
$ws3.Range("A3").Value = "Total filas: 27"

# --- Sheet 6203-6173: new row appended ---
Set-Row $ws3 32 "09:26:30" "11:14" "215C_LA PLATA" 108 "L6203"
